# "Generate Report for Handoff"
#
# The localization-status report is regenerated. The batch of files that
# were stamped with the previous handoff timestamp (zh-cn: 2016-03-10
# 16:23:23, de-de: 2016-03-10 16:23:27) - plus the two rows that had their
# own slightly later timestamps (zh-cn: 16:23:48, de-de: 16:23:53) - all
# now share one fresh "Latest Handoff Datetime" from this handoff run.

$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

foreach ($r in $rows) {
    $zh.Cells.Item($r, 4).Value = "2016-03-10 16:24:10"
}

foreach ($r in $rows) {
    $de.Cells.Item($r, 4).Value = "2016-03-10 16:24:15"
}
